# Weekly update: a new daily price record is inserted at row 81 for the
# "Feria Lagunitas de Puerto Montt - Zanahoria" subset, pushing the
# existing records (previously rows 81-184) down by one row to rows 82-185.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 81; this shifts rows 81:184 -> 82:185
# and Excel extends the sheet dimension/used range automatically.
$ws.Rows("81:81").Insert()

# Columns A,B,C,E,F,G,H,I,N,Q,R are constant for every record in this
# subset, so re-use the values that are now sitting in row 82 (the record
# that used to be row 81) to populate the new row.
$ws.Range("A81").Value = $ws.Range("A82").Value2
$ws.Range("B81").Value = $ws.Range("B82").Value()
$ws.Range("C81").Value = $ws.Range("C82").Value()
$ws.Range("E81").Value = $ws.Range("E82").Value2
$ws.Range("F81").Value = $ws.Range("F82").Value2
$ws.Range("G81").Value = $ws.Range("G82").Value()
$ws.Range("H81").Value = $ws.Range("H82").Value()
$ws.Range("I81").Value = $ws.Range("I82").Value()
$ws.Range("N81").Value = $ws.Range("N82").Value()
$ws.Range("Q81").Value = $ws.Range("Q82").Value2
$ws.Range("R81").Value = $ws.Range("R82").Value()

# New record-specific values for the newly inserted row.
$ws.Range("D81").Value = 44413
$ws.Range("J81").Value = 250
$ws.Range("K81").Value = 8000
$ws.Range("L81").Value = 8000
$ws.Range("M81").Value = 8000
$ws.Range("O81").Value = "Región de Ñuble"
$ws.Range("P81").Value = 400
